# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.405.28"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.622.02"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.210"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.12%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "4.200.47"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "606.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "70.480.29"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.631.59"
$ws.Range("E18").Value = "  +3.09%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("E26").Value = "  -5.45%  "
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.15%  "
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "0.0₃0882"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "3.909.26"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "519.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.20%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("E45").Value = "  +9.24%  "
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("E51").Value = "  +1.37%  "
